$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.216730117797852
$ws.Range("B1").Value = 2.162088632583618
$ws.Range("C1").Value = 4.053619384765625
$ws.Range("D1").Value = 3.160155296325684
$ws.Range("E1").Value = 1.096912264823914
